# Rename CLOCK and RESET columns on the SAFETY.PARITY sheet.
# - Insert two new columns (HSR ID, SM ID) before the old "IP NAME" column.
# - Rename "IP CLOCK NAME" -> "CLOCK" and "IP RESET NAME" -> "RESET".
# - Refresh the sheet view (zoom / top-left cell) to match the other sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SAFETY.PARITY")

# Insert two blank columns at B:C, shifting everything from the old B onward
# two columns to the right (old B -> D, old C -> E, ...).
$ws.Range("B:C").EntireColumn.Insert()

# New header cells for the inserted columns.
$ws.Range("B2").Value = "HSR ID"
$ws.Range("C2").Value = "SM ID"
$ws.Range("B2:C2").WrapText = $false

# The old "IP CLOCK NAME" / "IP RESET NAME" headers now live at E2 / F2 after
# the column insert - rename them to match the SAFETY.DCLS naming.
$ws.Range("E2").Value = "CLOCK"
$ws.Range("F2").Value = "RESET"

# Match the view settings used elsewhere in the workbook.
$ws.Application.ActiveWindow.Zoom = 85
$ws.Range("B1").Select()
